$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 33904
$ws.Range("D2").Value = 49175127
$ws.Range("C3").Value = 83038
$ws.Range("D3").Value = 121991959
$ws.Range("C4").Value = 28414
$ws.Range("D4").Value = 42143697
$ws.Range("C5").Value = 7743
$ws.Range("D5").Value = 11521714
$ws.Range("C6").Value = 1587
$ws.Range("D6").Value = 2365128
$ws.Range("C7").Value = 112
$ws.Range("D7").Value = 163093
$ws.Range("C10").Value = 36925
$ws.Range("D10").Value = 50396966
$ws.Range("C11").Value = 8743
$ws.Range("D11").Value = 12672550
$ws.Range("C12").Value = 23957
$ws.Range("D12").Value = 35204496
$ws.Range("C13").Value = 7606
$ws.Range("D13").Value = 11308298
$ws.Range("C14").Value = 1906
$ws.Range("D14").Value = 2835211
$ws.Range("C15").Value = 337
$ws.Range("D15").Value = 495964
$ws.Range("C18").Value = 9205
$ws.Range("D18").Value = 12264409
$ws.Range("C19").Value = 12253
$ws.Range("D19").Value = 17724167
$ws.Range("C20").Value = 29306
$ws.Range("D20").Value = 43097864
$ws.Range("C21").Value = 9494
$ws.Range("D21").Value = 14128387
$ws.Range("C22").Value = 2369
$ws.Range("D22").Value = 3529219
$ws.Range("C23").Value = 386
$ws.Range("D23").Value = 576342
$ws.Range("C25").Value = 10599
$ws.Range("D25").Value = 14279971
$ws.Range("C26").Value = 6966
$ws.Range("D26").Value = 10109693
$ws.Range("C27").Value = 20744
$ws.Range("D27").Value = 30500419
$ws.Range("C28").Value = 7206
$ws.Range("D28").Value = 10729039
$ws.Range("C29").Value = 1762
$ws.Range("D29").Value = 2636482
$ws.Range("C30").Value = 270
$ws.Range("D30").Value = 402415
$ws.Range("C32").Value = 7497
$ws.Range("D32").Value = 9970944
$ws.Range("C33").Value = 2699
$ws.Range("D33").Value = 3897059
$ws.Range("C34").Value = 6809
$ws.Range("D34").Value = 9957833
$ws.Range("C35").Value = 2760
$ws.Range("D35").Value = 4086732
$ws.Range("C36").Value = 730
$ws.Range("D36").Value = 1087763
$ws.Range("C37").Value = 130
$ws.Range("D37").Value = 194304
$ws.Range("C39").Value = 2020
$ws.Range("D39").Value = 2726045
$ws.Range("C40").Value = 15551
$ws.Range("D40").Value = 22539995
$ws.Range("C41").Value = 46822
$ws.Range("D41").Value = 68741146
$ws.Range("C42").Value = 17499
$ws.Range("D42").Value = 26005445
$ws.Range("C43").Value = 5011
$ws.Range("D43").Value = 7469155
$ws.Range("C44").Value = 923
$ws.Range("D44").Value = 1379950
$ws.Range("C48").Value = 14954
$ws.Range("D48").Value = 20042130
$ws.Range("C49").Value = 1599
$ws.Range("D49").Value = 2319448
$ws.Range("C50").Value = 5654
$ws.Range("D50").Value = 8326804
$ws.Range("C51").Value = 1993
$ws.Range("D51").Value = 2977957
$ws.Range("C52").Value = 644
$ws.Range("D52").Value = 962026
$ws.Range("C55").Value = 4731
$ws.Range("D55").Value = 6555057
$ws.Range("C56").Value = 609
$ws.Range("D56").Value = 891985
$ws.Range("C57").Value = 1559
$ws.Range("D57").Value = 2309923
$ws.Range("C58").Value = 607
$ws.Range("D58").Value = 906419
$ws.Range("C62").Value = 884
$ws.Range("D62").Value = 1250750
$ws.Range("C63").Value = 13915
$ws.Range("D63").Value = 20129217
$ws.Range("C64").Value = 41343
$ws.Range("D64").Value = 60605710
$ws.Range("C65").Value = 14537
$ws.Range("D65").Value = 21631516
$ws.Range("C66").Value = 4124
$ws.Range("D66").Value = 6148083
$ws.Range("C67").Value = 758
$ws.Range("D67").Value = 1130520
$ws.Range("C68").Value = 64
$ws.Range("D68").Value = 94689
$ws.Range("C70").Value = 13854
$ws.Range("D70").Value = 18365019
$ws.Range("C71").Value = 3946
$ws.Range("D71").Value = 5732057
$ws.Range("C72").Value = 10286
$ws.Range("D72").Value = 15124381
$ws.Range("C73").Value = 3519
$ws.Range("D73").Value = 5246526
$ws.Range("C74").Value = 1199
$ws.Range("D74").Value = 1793091
$ws.Range("C78").Value = 4570
$ws.Range("D78").Value = 6170443
$ws.Range("C79").Value = 1313
$ws.Range("D79").Value = 1897497
$ws.Range("C80").Value = 4381
$ws.Range("D80").Value = 6455062
$ws.Range("C81").Value = 1706
$ws.Range("D81").Value = 2546932
$ws.Range("C82").Value = 599
$ws.Range("D82").Value = 897641
$ws.Range("C83").Value = 141
$ws.Range("D83").Value = 211069
$ws.Range("C86").Value = 2839
$ws.Range("D86").Value = 3778957
$ws.Range("C87").Value = 488
$ws.Range("D87").Value = 728159
$ws.Range("C90").Value = 32
$ws.Range("D90").Value = 48000
$ws.Range("C92").Value = 9774
$ws.Range("D92").Value = 14229434
$ws.Range("C93").Value = 27088
$ws.Range("D93").Value = 39867767
$ws.Range("C94").Value = 9059
$ws.Range("D94").Value = 13481612
$ws.Range("C95").Value = 2444
$ws.Range("D95").Value = 3647399
$ws.Range("C96").Value = 405
$ws.Range("D96").Value = 605249
$ws.Range("C99").Value = 8875
$ws.Range("D99").Value = 11786200
$ws.Range("C100").Value = 27511
$ws.Range("D100").Value = 39754237
$ws.Range("C101").Value = 61094
$ws.Range("D101").Value = 89559023
$ws.Range("C102").Value = 19749
$ws.Range("D102").Value = 29384531
$ws.Range("C103").Value = 5445
$ws.Range("D103").Value = 8118982
$ws.Range("C104").Value = 938
$ws.Range("D104").Value = 1401893
$ws.Range("C105").Value = 51
$ws.Range("D105").Value = 73358
$ws.Range("C108").Value = 23389
$ws.Range("D108").Value = 31398223
$ws.Range("C109").Value = 31996
$ws.Range("D109").Value = 46269732
$ws.Range("C110").Value = 69701
$ws.Range("D110").Value = 102117634
$ws.Range("C111").Value = 21791
$ws.Range("D111").Value = 32369780
$ws.Range("C112").Value = 5704
$ws.Range("D112").Value = 8486863
$ws.Range("C113").Value = 977
$ws.Range("D113").Value = 1457677
$ws.Range("C114").Value = 47
$ws.Range("D114").Value = 68728
$ws.Range("C115").Value = 15
$ws.Range("D115").Value = 22500
$ws.Range("C117").Value = 28197
$ws.Range("D117").Value = 37636643
$ws.Range("C118").Value = 12045
$ws.Range("D118").Value = 17474245
$ws.Range("C119").Value = 30028
$ws.Range("D119").Value = 44165716
$ws.Range("C120").Value = 10650
$ws.Range("D120").Value = 15837716
$ws.Range("C121").Value = 2674
$ws.Range("D121").Value = 3989370
$ws.Range("C122").Value = 394
$ws.Range("D122").Value = 584990
$ws.Range("C125").Value = 9874
$ws.Range("D125").Value = 13268415
$ws.Range("C126").Value = 31280
$ws.Range("D126").Value = 45235774
$ws.Range("C127").Value = 73956
$ws.Range("D127").Value = 108495538
$ws.Range("C128").Value = 22267
$ws.Range("D128").Value = 33125553
$ws.Range("C129").Value = 5712
$ws.Range("D129").Value = 8529097
$ws.Range("C130").Value = 1173
$ws.Range("D130").Value = 1749220
$ws.Range("C133").Value = 26275
$ws.Range("D133").Value = 35684768
$ws.Range("C134").Value = 43838
$ws.Range("D134").Value = 63898294
$ws.Range("C135").Value = 127780
$ws.Range("D135").Value = 188535363
$ws.Range("C136").Value = 56098
$ws.Range("D136").Value = 83654110
$ws.Range("C137").Value = 17705
$ws.Range("D137").Value = 26462934
$ws.Range("C138").Value = 3756
$ws.Range("D138").Value = 5613773
$ws.Range("C139").Value = 201
$ws.Range("D139").Value = 299355
$ws.Range("C140").Value = 18
$ws.Range("D140").Value = 25905
$ws.Range("C145").Value = 43448
$ws.Range("D145").Value = 59651182
